$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 177-178, shifting existing rows 177-202 down to 179-204
$ws.Rows("177:178").Insert()

# Row 177
$ws.Range("A177").Value = 8
$ws.Range("B177").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C177").Value = 'Coquimbo'
$ws.Range("D177").Value = 44505
$ws.Range("E177").Value = 4
$ws.Range("F177").Value = 100112032
$ws.Range("G177").Value = 'Zapallo italiano'
$ws.Range("H177").Value = 'Bola 8'
$ws.Range("I177").Value = 'Primera'
$ws.Range("J177").Value = 400
$ws.Range("K177").Value = 7000
$ws.Range("L177").Value = 8000
$ws.Range("M177").Value = 7500
$ws.Range("N177").Value = '$/caja 60 unidades'
$ws.Range("O177").Value = 'Región de Arica y Parinacota'
$ws.Range("P177").Value = 125
$ws.Range("Q177").Value = 60
$ws.Range("R177").Value = 'Hortaliza'

# Row 178
$ws.Range("A178").Value = 8
$ws.Range("B178").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C178").Value = 'Coquimbo'
$ws.Range("D178").Value = 44505
$ws.Range("E178").Value = 4
$ws.Range("F178").Value = 100112032
$ws.Range("G178").Value = 'Zapallo italiano'
$ws.Range("H178").Value = 'Sin especificar'
$ws.Range("I178").Value = 'Primera'
$ws.Range("J178").Value = 520
$ws.Range("K178").Value = 10000
$ws.Range("L178").Value = 11000
$ws.Range("M178").Value = 10500
$ws.Range("N178").Value = '$/caja 70 unidades'
$ws.Range("O178").Value = 'Provincia de Limarí'
$ws.Range("P178").Value = 150
$ws.Range("Q178").Value = 70
$ws.Range("R178").Value = 'Hortaliza'
